$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source dataset dropped two samples (rows originally holding "RM 232"
# and "SC 92"); delete those two rows first so everything below shifts up,
# matching the new A1:F33 layout.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Final values for the data block (A2:F33) after the row removals and the
# additional missing-data corrections called out in the diff.
$data = @(
    @('RM 2', -19.7, 14.9, -13.5, $null, 18.03),
    @('RM 8', -19.7, 11.2, -14.2, -5.7, 17.64),
    @('RM 9', $null, 11, -15.4, -6.4, $null),
    @('RM 14', -19.5, 12.3, -14.4, -5, 17.66),
    @('RM 21', -19.8, 15.1, -14.2, -5.7, 16.43),
    @('RM 32', $null, 15, -13.8, -7.1, 17.24),
    @('RM 38', -19.9, $null, -13.9, -6.6, 17.05),
    @('RM 42', -20.6, 10.5, -14.5, -6.8, 17.26),
    @('RM 52 a', -19.8, 11.5, -14.7, -6.1, 16.43),
    @('RM 58', -20.8, 11.4, -15.5, -7.9, 17.65),
    @('RM 81', -19.9, 12.5, -14.1, -5.3, $null),
    @('RM 88', -19.9, 12.5, -13.9, -5.3, 17.1),
    @('RM 90', -19.6, 14.4, -13.1, -5.4, 17.76),
    @('RM 95', -19.1, 12.5, -15.2, -8.4, 16.2),
    @('RM 103', $null, 13.5, -15.3, -6.9, 17.34),
    @('RM 116', -19.4, 11.2, -14.7, -7.3, 17.78),
    @('RM 120', -19.6, 11.5, -15.2, -8.5, 18.35),
    @('RM 125', -20.6, 13.2, -15.5, -6.5, 17.81),
    @('RM 134', -19, $null, -14, -7.2, $null),
    @('RM 135', -18.9, 12.7, -14.3, -8.7, 16.58),
    @('RM 138', -19.3, 13.5, -15.4, -6.1, 16.81),
    @('RM 140', -19.5, 12.2, -13.9, -7, 16.48),
    @('RM 142a', -17.7, 12.7, -13.9, -8.1, 16.78),
    @('RM 145', -19.5, $null, -15.5, -7.1, 16.6),
    @('SC 5', -20.2, 10.8, -13.8, -5, $null),
    @('SC 101', -20.4, 10, -14.6, -10, $null),
    @('SC 105', $null, 11.1, -13.7, -5.9, $null),
    @('SC 119', $null, 11.2, -13, -6.8, 18.06),
    @('SC 120', -19.7, 11.4, -13.6, -5.7, 16.89),
    @('SC 132', -18.8, 15.3, -13.7, -8.1, 17.18),
    @('SC 193', $null, 10.5, -14.7, -6.4, 17.39),
    @('SC 232', -19.5, 10.4, -14.1, -10.7, 17.53)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
